$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.812.32"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.342.89"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +6.43%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.51"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.71"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.760.90"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.753.73"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000134"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.344.34"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.67"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.29"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "328.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "62.79"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.62%  "
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("E27").Value = "  -5.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.75"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.98"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0734"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.32"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("B33").Value = "SuiNetwork"
$ws.Range("C33").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.01"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.17"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.46%  "
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.60"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.06"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "141.37"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.80%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.64"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "287.35"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0947"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0511"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.09"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.566"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.954"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.98%  "
